$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for rows 2-176.
# All of them move forward by one day (45184 -> 45185).
for ($r = 2; $r -le 176; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
